$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This change records a new handoff report entry for file
#   1263bb52-8566-4b30-95f5-ea042f61c688.md
# which becomes the newest row (row 2) on every sheet, pushing the previously
# newest entry (994e2a69-aee8-4e51-aeb9-1bf4ce473faa.md) down to row 3.
# Row 3 is an exact duplicate of the old row 2 (same values/format), so the
# cheapest & most faithful way to reproduce it is: copy row 2 -> row 3, fix
# up the hyperlinks, then overwrite the handful of cells on row 2 that are
# actually different for the new file.
# ---------------------------------------------------------------------------

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b08209476cd98422d97d0d7c225dc9d7adf6eddb/e2e/"

$oldFile = "994e2a69-aee8-4e51-aeb9-1bf4ce473faa.md"
$newFile = "1263bb52-8566-4b30-95f5-ea042f61c688.md"

$oldDisplay = "e2e\994e2a69-aee8-4e51-aeb9-1bf4ce473faa.md"
$newDisplay = "e2e\1263bb52-8566-4b30-95f5-ea042f61c688.md"

# ===========================================================================
# Sheet "Overview" (columns A:G, hyperlink lives in column B)
# ===========================================================================
$ws = $wb.Worksheets.Item("Overview")

# Duplicate row 2 -> row 3 (keeps all values + formatting for the old file)
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()

# The old hyperlink object is still anchored on B2 (text now lives on B3) -
# drop it and recreate clean hyperlinks on B2 (new file) / B3 (old file).
$ws.Range("B2").Hyperlinks.Delete()

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = $newDisplay
$ws.Range("G2").Value = "2016-08-27 04:38:33"
$ws.Hyperlinks.Add($ws.Range("B2"), ($githubBase + $newFile), "", "", $newDisplay)

$ws.Range("B3").Value = $oldDisplay
$ws.Hyperlinks.Add($ws.Range("B3"), ($githubBase + $oldFile), "", "", $oldDisplay)

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G3"))

# ===========================================================================
# Sheet "zh-cn" (columns A:P, hyperlink lives in column A)
# ===========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()

$ws.Range("A2").Hyperlinks.Delete()

$ws.Range("A2").Value = $newFile
$ws.Range("G2").Value = "1263bb52-8566-4b30-95f5-ea042f61c688.ff9ca852a45243b004d5194abfd355b39d6354f5.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-27 04:38:29"
$ws.Hyperlinks.Add($ws.Range("A2"), ($githubBase + $newFile), "", "", $newFile)

$ws.Range("A3").Value = $oldFile
$ws.Hyperlinks.Add($ws.Range("A3"), ($githubBase + $oldFile), "", "", $oldFile)

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))

# ===========================================================================
# Sheet "de-de" (columns A:P, hyperlink lives in column A)
# ===========================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()

$ws.Range("A2").Hyperlinks.Delete()

$ws.Range("A2").Value = $newFile
$ws.Range("G2").Value = "1263bb52-8566-4b30-95f5-ea042f61c688.ff9ca852a45243b004d5194abfd355b39d6354f5.de-de.xlf"
$ws.Range("H2").Value = "2016-08-27 04:38:33"
$ws.Hyperlinks.Add($ws.Range("A2"), ($githubBase + $newFile), "", "", $newFile)

$ws.Range("A3").Value = $oldFile
$ws.Hyperlinks.Add($ws.Range("A3"), ($githubBase + $oldFile), "", "", $oldFile)

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))
